$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Fix comma separators to periods in contractor name strings ---
$ws.Range("E16").Value = "MARSICO GUILLERMO MIGUEL. MARSICO JUAN EDUARDO"
$ws.Range("E84").Value = "MARSICO GUILLERMO MIGUEL. MARSICO JUAN EDUARDO"
$ws.Range("E20").Value = "IZAGUIRRE CARLOS MARIA. MOREND MARIA ELENA Y MOREND MARIA TERESA"
$ws.Range("F20").Value = "IZAGUIRRE CARLOS MARIA. MOREND MARIA ELENA Y MOREND MARIA TERESA"
$ws.Range("E43").Value = "IZAGUIRRE CARLOS MARIA. MOREND MARIA ELENA Y MOREND MARIA TERESA"
$ws.Range("F43").Value = "IZAGUIRRE CARLOS MARIA. MOREND MARIA ELENA Y MOREND MARIA TERESA"
$ws.Range("E42").Value = "FERNANDEZ MARIO H. GALLICET OSCAR M"

# --- Fix decimal number formatting (comma decimal/dot thousands -> dot decimal, no thousands sep) ---
# Each amount is written as a formula returning a text literal, then the whole range is
# converted to static values in a single paste-special pass. This keeps the cell a plain
# text (shared-string) cell using its original default style, instead of Excel silently
# coercing the numeric-looking string into a Number value (which would also re-style it).
$ws.Range("H2").Formula = "=""240.00"""
$ws.Range("H3").Formula = "=""2100.00"""
$ws.Range("H4").Formula = "=""12524.40"""
$ws.Range("H5").Formula = "=""89934.52"""
$ws.Range("H6").Formula = "=""720.00"""
$ws.Range("H7").Formula = "=""5122.49"""
$ws.Range("H8").Formula = "=""7540.47"""
$ws.Range("H9").Formula = "=""591.00"""
$ws.Range("H10").Formula = "=""10034.79"""
$ws.Range("H11").Formula = "=""14.10"""
$ws.Range("H12").Formula = "=""143.96"""
$ws.Range("H13").Formula = "=""140.00"""
$ws.Range("H14").Formula = "=""19563.86"""
$ws.Range("H15").Formula = "=""101.60"""
$ws.Range("H16").Formula = "=""5190.00"""
$ws.Range("H17").Formula = "=""144.60"""
$ws.Range("H18").Formula = "=""11990.00"""
$ws.Range("H19").Formula = "=""2569.00"""
$ws.Range("H20").Formula = "=""519.20"""
$ws.Range("H21").Formula = "=""3564.00"""
$ws.Range("H22").Formula = "=""19362.78"""
$ws.Range("H23").Formula = "=""8012.00"""
$ws.Range("H24").Formula = "=""1005.00"""
$ws.Range("H25").Formula = "=""874.74"""
$ws.Range("H26").Formula = "=""84.50"""
$ws.Range("H27").Formula = "=""1128.00"""
$ws.Range("H28").Formula = "=""6.78"""
$ws.Range("H29").Formula = "=""691.00"""
$ws.Range("H30").Formula = "=""33542.88"""
$ws.Range("H31").Formula = "=""7756.64"""
$ws.Range("H32").Formula = "=""72.00"""
$ws.Range("H33").Formula = "=""2384.50"""
$ws.Range("H34").Formula = "=""588.00"""
$ws.Range("H35").Formula = "=""7982.32"""
$ws.Range("H36").Formula = "=""194.32"""
$ws.Range("H37").Formula = "=""15447.00"""
$ws.Range("H38").Formula = "=""890.00"""
$ws.Range("H39").Formula = "=""3760.00"""
$ws.Range("H40").Formula = "=""43077.24"""
$ws.Range("H41").Formula = "=""9397.50"""
$ws.Range("H42").Formula = "=""647.00"""
$ws.Range("H43").Formula = "=""478.38"""
$ws.Range("H44").Formula = "=""9030.00"""
$ws.Range("H45").Formula = "=""199.00"""
$ws.Range("H46").Formula = "=""4.38"""
$ws.Range("H47").Formula = "=""144870.00"""
$ws.Range("H48").Formula = "=""34000.00"""
$ws.Range("H49").Formula = "=""3.60"""
$ws.Range("H50").Formula = "=""70.96"""
$ws.Range("H51").Formula = "=""62.00"""
$ws.Range("H52").Formula = "=""10000.00"""
$ws.Range("H53").Formula = "=""26506.90"""
$ws.Range("H54").Formula = "=""2156.30"""
$ws.Range("H55").Formula = "=""1327.00"""
$ws.Range("H56").Formula = "=""36.00"""
$ws.Range("H57").Formula = "=""7550.00"""
$ws.Range("H58").Formula = "=""561.75"""
$ws.Range("H59").Formula = "=""466.00"""
$ws.Range("H60").Formula = "=""8889.00"""
$ws.Range("H61").Formula = "=""11116.00"""
$ws.Range("H62").Formula = "=""73500.00"""
$ws.Range("H63").Formula = "=""14052.00"""
$ws.Range("H64").Formula = "=""902.00"""
$ws.Range("H65").Formula = "=""2000.00"""
$ws.Range("H66").Formula = "=""4813.38"""
$ws.Range("H67").Formula = "=""500.00"""
$ws.Range("H68").Formula = "=""250.00"""
$ws.Range("H69").Formula = "=""1800.00"""
$ws.Range("H70").Formula = "=""4000.00"""
$ws.Range("H71").Formula = "=""24723.25"""
$ws.Range("H72").Formula = "=""500.00"""
$ws.Range("H73").Formula = "=""2850.00"""
$ws.Range("H74").Formula = "=""650.00"""
$ws.Range("H75").Formula = "=""250.00"""
$ws.Range("H76").Formula = "=""3920.00"""
$ws.Range("H77").Formula = "=""1815.00"""
$ws.Range("H78").Formula = "=""1500.00"""
$ws.Range("H79").Formula = "=""200.00"""
$ws.Range("H80").Formula = "=""240.00"""
$ws.Range("H81").Formula = "=""1000.00"""
$ws.Range("H82").Formula = "=""3642.90"""
$ws.Range("H83").Formula = "=""1560.90"""
$ws.Range("H84").Formula = "=""670.00"""
$ws.Range("H85").Formula = "=""4385.00"""
$ws.Range("H86").Formula = "=""2581.40"""
$ws.Range("H87").Formula = "=""2581.06"""
$ws.Range("H88").Formula = "=""1290.00"""
$ws.Range("H89").Formula = "=""511.53"""
$ws.Range("H90").Formula = "=""288.00"""
$ws.Range("H91").Formula = "=""6032.00"""
$ws.Range("H92").Formula = "=""930.00"""
$ws.Range("H93").Formula = "=""490.00"""
$ws.Range("H94").Formula = "=""117.90"""
$ws.Range("H95").Formula = "=""2452.50"""
$ws.Range("H96").Formula = "=""630.00"""
$ws.Range("H97").Formula = "=""740.40"""
$ws.Range("H98").Formula = "=""2119.64"""
$ws.Range("H99").Formula = "=""4438.97"""
$ws.Range("H100").Formula = "=""1380.01"""
$ws.Range("H101").Formula = "=""1900.00"""
$ws.Range("H102").Formula = "=""9780.00"""
$ws.Range("H103").Formula = "=""580.00"""
$ws.Range("H104").Formula = "=""2101.94"""
$ws.Range("H105").Formula = "=""1070.00"""
$ws.Range("H106").Formula = "=""116667.00"""
$ws.Range("H107").Formula = "=""15000.00"""
$ws.Range("H108").Formula = "=""116667.00"""
$ws.Range("H109").Formula = "=""116667.00"""
$ws.Range("H110").Formula = "=""33300.00"""
$ws.Range("H111").Formula = "=""49000.00"""
$ws.Range("H112").Formula = "=""78000.00"""
$ws.Range("H113").Formula = "=""1720.00"""
$ws.Range("H114").Formula = "=""651.00"""

$rng = $ws.Range("H2:H114")
$rng.Copy()
$rng.PasteSpecial(-4163)
$excel.CutCopyMode = 0

